$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Helper: find the Hyperlink object currently anchored to a given Range.
# ---------------------------------------------------------------------------
function Get-HyperlinkForRange($ws, $range) {
    $target = $range.Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            return $hl
        }
    }
    return $null
}

# Helper: copy the "hyperlink" visual style (underline + blue font) that the
# workbook already uses for its other link cells (e.g. column A / D).
function Set-LinkStyle($range) {
    $range.Font.Underline = 2          # xlUnderlineStyleSingle
    $range.Font.Color = 15570276       # RGB(100,149,237) == "FF6495ED"
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
}

# ===========================================================================
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used on the Overview sheet (B/C) and on each
#    language sheet's Status column (C).
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ===========================================================================
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File"
#    (G) for the handed-back rows, and refresh the handback datetime (H).
# ===========================================================================
foreach ($row in 2, 3) {
    $aCell = $wsZhCn.Range("A$row")
    $dCell = $wsZhCn.Range("D$row")
    $fCell = $wsZhCn.Range("F$row")
    $gCell = $wsZhCn.Range("G$row")

    $aLink = Get-HyperlinkForRange $wsZhCn $aCell
    $dLink = Get-HyperlinkForRange $wsZhCn $dCell

    $fCell.Value = $aCell.Text
    Set-LinkStyle $fCell
    if ($aLink -ne $null) {
        $wsZhCn.Hyperlinks.Add($fCell, $aLink.Address, "", "", $aCell.Text) | Out-Null
    }

    $gCell.Value = $dCell.Text
    Set-LinkStyle $gCell
    if ($dLink -ne $null) {
        $wsZhCn.Hyperlinks.Add($gCell, $dLink.Address, "", "", $dCell.Text) | Out-Null
    }
}

# Handback completed for zh-cn.
$wsZhCn.Range("H2").Value = "2016-03-20 20:52:41"
$wsZhCn.Range("H3").Value = "2016-03-20 20:52:41"

# ===========================================================================
# 3. de-de sheet: same shape of change, different timestamp / targets.
# ===========================================================================
foreach ($row in 2, 3) {
    $aCell = $wsDeDe.Range("A$row")
    $dCell = $wsDeDe.Range("D$row")
    $fCell = $wsDeDe.Range("F$row")
    $gCell = $wsDeDe.Range("G$row")

    $aLink = Get-HyperlinkForRange $wsDeDe $aCell
    $dLink = Get-HyperlinkForRange $wsDeDe $dCell

    $fCell.Value = $aCell.Text
    Set-LinkStyle $fCell
    if ($aLink -ne $null) {
        $wsDeDe.Hyperlinks.Add($fCell, $aLink.Address, "", "", $aCell.Text) | Out-Null
    }

    $gCell.Value = $dCell.Text
    Set-LinkStyle $gCell
    if ($dLink -ne $null) {
        $wsDeDe.Hyperlinks.Add($gCell, $dLink.Address, "", "", $dCell.Text) | Out-Null
    }
}

# Handback completed for de-de (later than zh-cn -> distinct timestamp).
$wsDeDe.Range("H2").Value = "2016-03-20 20:52:47"
$wsDeDe.Range("H3").Value = "2016-03-20 20:52:47"
